$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E24").Value = 3819.58
$ws.Range("C26").Value = 43.4
$ws.Range("E26").Value = 7878.947
